# --- Fix header labels on the existing sheets -------------------------
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet after the last sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy the header formatting (bold, centered, bordered) from the
# "Weekly Quantity" sheet and stamp it across the four header cells.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Copy the date number-format used for column A of "Weekly Quantity"
# down the 30 data rows of the new sheet.
$ws1.Range("A2").Copy()
$ws3.Range("A2:A31").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$poForecastData = @(
    @(44941.99999999999, 7, -0.2529461816881701, 13.45614474312756),
    @(44948.99999999999, 7, -0.2001062006925848, 13.61927558843259),
    @(45004.99999999999, 7, 0.9144775973765219, 14.55230068048852),
    @(45011.99999999999, 7, 0.716366537840976, 14.29901449555035),
    @(45018.99999999999, 7, 0.5300781986441457, 14.91001860702349),
    @(45025.99999999999, 8, 0.1240553079799721, 14.60075462958339),
    @(45032.99999999999, 8, 0.5374278546221154, 14.647784764534),
    @(45039.99999999999, 8, 0.8168153585000472, 14.80988313340106),
    @(45046.99999999999, 8, 0.2529234001344964, 14.76938307219815),
    @(45053.99999999999, 8, 1.074600885031107, 14.88804913918709),
    @(45060.99999999999, 8, 0.6411736049327722, 15.02632233187541),
    @(45067.99999999999, 8, 0.5090242335647437, 15.10156669788183),
    @(45074.99999999999, 8, 1.309119122552987, 15.52427186108896),
    @(45081.99999999999, 8, 1.172652390230297, 15.04468603118458),
    @(45515.99999999999, 13, 6.6216339657128, 20.54527982074802),
    @(45536.99999999999, 14, 6.865095179145324, 20.65951566594298),
    @(45543.99999999999, 14, 6.618366891886208, 20.89182330390041),
    @(45550.99999999999, 14, 6.775333095319407, 20.45356002553279),
    @(45557.99999999999, 14, 6.847339899714221, 20.90653876370589),
    @(45564.99999999999, 14, 6.404646037295858, 20.8794325776808),
    @(45571.99999999999, 14, 6.748802490744437, 21.28220637500882),
    @(45578.99999999999, 14, 7.42162876140272, 21.03509427243883),
    @(45585.99999999999, 14, 7.748919023086518, 21.20680417359095),
    @(45592.99999999999, 14, 7.22580349938072, 21.11637952867225),
    @(45599.99999999999, 14, 7.424454561447041, 21.54500872051289),
    @(45606.99999999999, 15, 7.356449824275005, 21.19601169522509),
    @(45613.99999999999, 15, 7.350448421849667, 21.57608314931372),
    @(45620.99999999999, 15, 7.337900283894497, 21.94828359320979),
    @(45627.99999999999, 15, 7.884594047080689, 21.99002851709307),
    @(45634.99999999999, 15, 7.939382827519296, 21.89584363765604)
)

for ($r = 0; $r -lt $poForecastData.Length; $r++) {
    $rowValues = $poForecastData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws3.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}

$ws3.Range("A1").Select()

Write-Output "PO Forecast sheet added with $($poForecastData.Length) data rows"
